$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) cells that hold numeric-looking text retain their
# original text representation (no auto-number conversion / precision loss).
$dCells = @("D2", "D3", "D4", "D5", "D7", "D8", "D9", "D11", "D12", "D13", "D14", "D15", "D16", "D17", "D18", "D19", "D20", "D21", "D24", "D25", "D27", "D28", "D29", "D30", "D31", "D32", "D34", "D35", "D36", "D37", "D38", "D39", "D40", "D41", "D42", "D43", "D44", "D45", "D46", "D48", "D50", "D51")
foreach ($cellRef in $dCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

$ws.Range("D2").Value = '24.710.01'
$ws.Range("E2").Value = '  +0.47%  '

$ws.Range("D3").Value = '1.694.10'
$ws.Range("E3").Value = '  -0.05%  '

$ws.Range("D4").Value = '1.003'
$ws.Range("E4").Value = '  +0.28%  '

$ws.Range("D5").Value = '316.92'
$ws.Range("E5").Value = '  +1.13%  '

$ws.Range("E6").Value = '  +0.27%  '

$ws.Range("D7").Value = '0.3955'
$ws.Range("E7").Value = '  +0.11%  '

$ws.Range("D8").Value = '0.4068'
$ws.Range("E8").Value = '  +0.85%  '

$ws.Range("D9").Value = '1.494'
$ws.Range("E9").Value = '  -1.81%  '

$ws.Range("E10").Value = '  +0.26%  '

$ws.Range("D11").Value = '53.04'
$ws.Range("E11").Value = '  -6.59%  '

$ws.Range("D12").Value = '0.08952'
$ws.Range("E12").Value = '  +2.03%  '

$ws.Range("D13").Value = '7.273'
$ws.Range("E13").Value = '  -0.84%  '

$ws.Range("D14").Value = '23.60'
$ws.Range("E14").Value = '  +2.26%  '

$ws.Range("D15").Value = '8.051'
$ws.Range("E15").Value = '  +5.63%  '

$ws.Range("D16").Value = '0.00001324'
$ws.Range("E16").Value = '  +0.39%  '

$ws.Range("D17").Value = '1.695.41'
$ws.Range("E17").Value = '  +0.20%  '

$ws.Range("D18").Value = '100.01'
$ws.Range("E18").Value = '  -0.36%  '

$ws.Range("D19").Value = '0.07040'
$ws.Range("E19").Value = '  -0.24%  '

$ws.Range("D20").Value = '19.64'
$ws.Range("E20").Value = '  +0.94%  '

$ws.Range("D21").Value = '6.997'
$ws.Range("E21").Value = '  +4.29%  '

$ws.Range("E22").Value = '  +0.20%  '

$ws.Range("E23").Value = '  +1.39%  '

$ws.Range("D24").Value = '24.693.47'
$ws.Range("E24").Value = '  +0.50%  '

$ws.Range("D25").Value = '3.285'
$ws.Range("E25").Value = '  +8.39%  '

$ws.Range("E26").Value = '  +2.11%  '

$ws.Range("D27").Value = '22.74'
$ws.Range("E27").Value = '  +1.74%  '

$ws.Range("D28").Value = '162.18'
$ws.Range("E28").Value = '  +1.53%  '

$ws.Range("D29").Value = '136.10'
$ws.Range("E29").Value = '  +1.92%  '

$ws.Range("D30").Value = '5.199'
$ws.Range("E30").Value = '  +0.41%  '

$ws.Range("D31").Value = '7.514'
$ws.Range("E31").Value = '  -1.46%  '

$ws.Range("D32").Value = '0.08658'
$ws.Range("E32").Value = '  +0.83%  '

$ws.Range("E33").Value = '  -3.21%  '

$ws.Range("D34").Value = '7.073'
$ws.Range("E34").Value = '  -4.10%  '

$ws.Range("D35").Value = '11.44'
$ws.Range("E35").Value = '  +4.02%  '

$ws.Range("D36").Value = '0.2742'
$ws.Range("E36").Value = '  +0.98%  '

$ws.Range("B37").Value = 'Aptos'
$ws.Range("C37").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D37").Value = '14.48'
$ws.Range("E37").Value = '  -1.64%  '

$ws.Range("B38").Value = 'WEMIXTOKEN'
$ws.Range("C38").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D38").Value = '1.881'
$ws.Range("E38").Value = '  -4.19%  '

$ws.Range("D39").Value = '0.09256'
$ws.Range("E39").Value = '  +2.94%  '

$ws.Range("D40").Value = '0.02727'
$ws.Range("E40").Value = '  -0.51%  '

$ws.Range("D41").Value = '1.473'
$ws.Range("E41").Value = '  +0.33%  '

$ws.Range("D42").Value = '0.7676'
$ws.Range("E42").Value = '  +0.65%  '

$ws.Range("D43").Value = '16.20'
$ws.Range("E43").Value = '  +5.25%  '

$ws.Range("D44").Value = '2.596'
$ws.Range("E44").Value = '  +6.03%  '

$ws.Range("D45").Value = '0.7172'
$ws.Range("E45").Value = '  +0.15%  '

$ws.Range("D46").Value = '4.215'
$ws.Range("E46").Value = '  +1.03%  '

$ws.Range("E47").Value = '  +0.28%  '

$ws.Range("D48").Value = '140.38'
$ws.Range("E48").Value = '  +0.00%  '

$ws.Range("E49").Value = '  -0.20%  '

$ws.Range("D50").Value = '91.03'
$ws.Range("E50").Value = '  +5.52%  '

$ws.Range("D51").Value = '0.07981'
$ws.Range("E51").Value = '  -0.21%  '
